# Add 2022-Q1 data
# -----------------
# The workbook has quarterly sheets (2021-Q1 .. 2021-Q4) plus a "总计"
# (totals) summary sheet. This script:
#   1. Creates a new "2022-Q1" sheet (positioned right before "总计"),
#      modelled on the most recent quarterly sheet so it inherits the
#      same look & feel (bold/bordered header row + index column).
#   2. Fills it with the 2022-Q1 fund holdings data.
#   3. Inserts a new top data row in "总计" for 2022-Q1 and renumbers/
#      shifts the existing rows down, preserving styles.

$wb = $excel.ActiveWorkbook

# Helper: force a numeric-looking string to be stored as TEXT (t="inlineStr"/t="s"),
# then strip the resulting style back to the sheet's plain (unstyled) data cells,
# matching how the other text cells in these sheets look (no explicit style).
function Set-TextCell($ws, $addr, $val) {
    $ws.Range($addr).Value = "'" + $val
    $ws.Range($addr).Style = "Normal"
}

# ---------------------------------------------------------------------------
# 1) Create the "2022-Q1" worksheet by cloning "2021-Q4" (keeps header/index
#    column styling identical) and place it immediately before "总计".
# ---------------------------------------------------------------------------
$src = $wb.Worksheets.Item("2021-Q4")
$totalSheet = $wb.Worksheets.Item("总计")
$src.Copy($totalSheet)
$q1ws = $wb.Worksheets.Item("2021-Q4 (2)")
$q1ws.Name = "2022-Q1"

# The cloned sheet has 8 data rows (rows 2-8); the new data needs 9 (rows 2-9),
# so stamp row 9's index cell (column A) with the same style as the others.
$q1ws.Range("A2").Copy()
$q1ws.Range("A9").PasteSpecial(-4122)
$q1ws.Application.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 2) Populate "2022-Q1" header + data.
# ---------------------------------------------------------------------------
$q1ws.Range("B1").Value = "基金代码"
$q1ws.Range("C1").Value = "基金名称"
$q1ws.Range("D1").Value = "基金规模"
$q1ws.Range("E1").Value = "股票总仓位"
$q1ws.Range("F1").Value = "仓位占比"
$q1ws.Range("G1").Value = "持有市值(亿元)"
$q1ws.Range("H1").Value = "仓位排名"

$q1ws.Range("A2").Value = 0
Set-TextCell $q1ws "B2" "516150"
Set-TextCell $q1ws "C2" "嘉实中证稀土产业ETF"
Set-TextCell $q1ws "D2" "25.17"
Set-TextCell $q1ws "E2" "99.75"
Set-TextCell $q1ws "F2" "3.90"
Set-TextCell $q1ws "G2" "0.9816"
$q1ws.Range("H2").Value = 9

$q1ws.Range("A3").Value = 1
Set-TextCell $q1ws "B3" "516780"
Set-TextCell $q1ws "C3" "华泰柏瑞中证稀土产业ETF"
Set-TextCell $q1ws "D3" "11.06"
Set-TextCell $q1ws "E3" "98.70"
Set-TextCell $q1ws "F3" "3.92"
Set-TextCell $q1ws "G3" "0.4336"
$q1ws.Range("H3").Value = 9

$q1ws.Range("A4").Value = 2
Set-TextCell $q1ws "B4" "159715"
Set-TextCell $q1ws "C4" "易方达中证稀土产业ETF"
Set-TextCell $q1ws "D4" "3.42"
Set-TextCell $q1ws "E4" "99.06"
Set-TextCell $q1ws "F4" "3.85"
Set-TextCell $q1ws "G4" "0.1317"
$q1ws.Range("H4").Value = 9

$q1ws.Range("A5").Value = 3
Set-TextCell $q1ws "B5" "159713"
Set-TextCell $q1ws "C5" "富国中证稀土产业交易型开放式指数证券投资基金"
Set-TextCell $q1ws "D5" "3.26"
Set-TextCell $q1ws "E5" "99.26"
Set-TextCell $q1ws "F5" "3.89"
Set-TextCell $q1ws "G5" "0.1268"
$q1ws.Range("H5").Value = 9

$q1ws.Range("A6").Value = 4
Set-TextCell $q1ws "B6" "005947"
Set-TextCell $q1ws "C6" "德邦民裕进取量化精选灵活配置混合A"
Set-TextCell $q1ws "D6" "0.53"
Set-TextCell $q1ws "E6" "94.44"
Set-TextCell $q1ws "F6" "6.83"
Set-TextCell $q1ws "G6" "0.0362"
$q1ws.Range("H6").Value = 5

$q1ws.Range("A7").Value = 5
Set-TextCell $q1ws "B7" "014331"
Set-TextCell $q1ws "C7" "华泰柏瑞中证稀土产业ETF联接A"
Set-TextCell $q1ws "D7" "0.86"
Set-TextCell $q1ws "E7" "24.22"
Set-TextCell $q1ws "F7" "1.08"
Set-TextCell $q1ws "G7" "0.0093"
$q1ws.Range("H7").Value = 9

$q1ws.Range("A8").Value = 6
Set-TextCell $q1ws "B8" "014332"
Set-TextCell $q1ws "C8" "华泰柏瑞中证稀土产业ETF联接C"
Set-TextCell $q1ws "D8" "0.70"
Set-TextCell $q1ws "E8" "24.22"
Set-TextCell $q1ws "F8" "1.08"
Set-TextCell $q1ws "G8" "0.0076"
$q1ws.Range("H8").Value = 9

$q1ws.Range("A9").Value = 7
Set-TextCell $q1ws "B9" "005948"
Set-TextCell $q1ws "C9" "德邦民裕进取量化精选灵活配置混合C"
Set-TextCell $q1ws "D9" "0.09"
Set-TextCell $q1ws "E9" "94.44"
Set-TextCell $q1ws "F9" "6.83"
Set-TextCell $q1ws "G9" "0.0061"
$q1ws.Range("H9").Value = 5

# ---------------------------------------------------------------------------
# 3) Update "总计": insert the new 2022-Q1 summary row at the top of the
#    data (row 2) and shift the previously-existing rows down by one,
#    renumbering the index column (A) as we go.
# ---------------------------------------------------------------------------
$tws = $wb.Worksheets.Item("总计")

# Stamp row 6's index cell with the same style used by the other index cells.
$tws.Range("A2").Copy()
$tws.Range("A6").PasteSpecial(-4122)
$tws.Application.CutCopyMode = $false

$tws.Range("A6").Value = 4
$tws.Range("B6").Value = $tws.Range("B5").Value()
$tws.Range("C6").Value = $tws.Range("C5").Value()
$tws.Range("D6").Value = $tws.Range("D5").Value()

$tws.Range("A5").Value = 3
$tws.Range("B5").Value = $tws.Range("B4").Value()
$tws.Range("C5").Value = $tws.Range("C4").Value()
$tws.Range("D5").Value = $tws.Range("D4").Value()

$tws.Range("A4").Value = 2
$tws.Range("B4").Value = $tws.Range("B3").Value()
$tws.Range("C4").Value = $tws.Range("C3").Value()
$tws.Range("D4").Value = $tws.Range("D3").Value()

$tws.Range("A3").Value = 1
$tws.Range("B3").Value = $tws.Range("B2").Value()
$tws.Range("C3").Value = $tws.Range("C2").Value()
$tws.Range("D3").Value = $tws.Range("D2").Value()

$tws.Range("A2").Value = 0
$tws.Range("B2").Value = "2022-Q1"
$tws.Range("C2").Value = 8
$tws.Range("D2").Value = 1.73

# Restore the originally-active sheet/selection (the "2021-Q4" -> "2022-Q1"
# copy operation leaves the new sheet activated as a side effect).
$wb.Worksheets.Item("2021-Q1").Activate()
